$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = "Administrator, Dr. Amira Sobhy, Dr. Veronia Rafat, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad"
$ws.Range("G3").Value = "Dr. Eman Tantawi, Administrator, Dr. Majorelle Magdy, Dr. Veronia Rafat, Dr. Hend Mahmoud, Dr. Asmaa Reda"
$ws.Range("G4").Value = "Dr. Eman Tantawi, Dr. Majorelle Magdy, Dr. Gehan Adel, Dr. Servinaz Sayed Mohammad, Dr. Hend Mahmoud, Dr. Asmaa Reda"
$ws.Range("G5").Value = "Dr. Eman Tantawi, Dr. Asmaa Reda, Dr. Amira Sobhy, Dr. Veronia Rafat"
$ws.Range("G6").Value = "Dr. Alshimaa Atef, Dr. Manar Montaser, Dr. Majorelle Magdy, Dr. Mohammad El-Tanany, Dr. Menna tuâ€™Allah Medhat"
$ws.Range("G7").Value = "Dr. Lamiaa Ossama, Dr. Fatma Elhady, Dr. Menna tu'Alllah Mohammad, Dr. Kerelos Zareef, Dr. Amera Ahmad Saad, Dr. Nada Mohammad, Dr. Abeer Ragab"
$ws.Range("G11").Value = "Dr. Aya Saeed, Dr. Safa Hany, Dr. Amal Awwad"
$ws.Range("G12").Value = "Dr. Dina Adel, Dr. Madeha Saeed, Dr. Eman M. Abo-Sakaya, Dr. Yasmeena Fattoh, Dr. Marina Youhanna, Dr. Amira Ibrahim"
$ws.Range("G13").Value = "Dr. Esraa Mostafa, Dr. Yasmeena Fattoh, Dr. Amira Ibrahim"
$ws.Range("G19").Value = "Dr. Mariam Toma Gerges, Dr. Rania Ahmad Youssef"
$ws.Range("G24").Value = "Dr. Youstina Gamil, Dr. Sarah Mahdy"
$ws.Range("G25").Value = "Menna tuâ€™Allah Gamil, Dr. Nouran Mahmoud"
$ws.Range("G27").Value = "Dr. Nourham Mostafa, Dr. Hana Amr"
$ws.Range("G28").Value = "Dr. Aya Emad, Dr. Maryam Ashraf"
$ws.Range("G30").Value = "Dr. Yassmen Ahmad, Dr. Aya Hanafy, Dr. Shorok Mohammad, Dr. Wafaa Ebida"
